$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits after the
#    ObjectIdentifier description paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Week 4 table: fill in "Project Status: " with "red".
# ------------------------------------------------------------------
$week4 = $d.Tables.Item(4)
$statusCell = $week4.Cell(1, 2)
$statusCell.Range.Text = "Project Status: red"
$statusCell.Range.Font.Size = 10
$statusCell.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# 3. Week 4 table: fill in the (empty) Comments cell with the
#    weekly comment text, and re-add the "_GoBack" bookmark right
#    after the new text (this is where Word leaves it after the
#    last edit made in a session).
# ------------------------------------------------------------------
$commentsCell = $week4.Cell(3, 2)
$commentsCell.Range.Text = "I did not put enough time into my project this week. Will do better next week."
$commentsCell.Range.Font.Size = 10
$commentsCell.Range.Font.SizeBi = 10

$endOfComment = $commentsCell.Range
$endOfComment.End = $endOfComment.End - 1
$endOfComment.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endOfComment) | Out-Null
